$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.030.90'
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").Value = '2.351.52'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E5").Value = '  +0.99%  '
$ws.Range("D6").Value = "'238.87"
$ws.Range("E6").Value = '  +1.53%  '
$ws.Range("D7").Value = "'73.97"
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("D9").Value = "'0.586"
$ws.Range("E9").Value = '  +8.09%  '
$ws.Range("E10").Value = '  +1.60%  '
$ws.Range("D11").Value = "'57.17"
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("D12").Value = "'32.19"
$ws.Range("E12").Value = '  +13.74%  '
$ws.Range("E13").Value = '  +1.02%  '
$ws.Range("D14").Value = "'7.20"
$ws.Range("E14").Value = '  +7.57%  '
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("E16").Value = '  -1.05%  '
$ws.Range("D17").Value = "'0.899"
$ws.Range("E17").Value = '  +1.01%  '
$ws.Range("D18").Value = '2.373.80'
$ws.Range("E18").Value = '  +0.87%  '
$ws.Range("D19").Value = '43.903.36'
$ws.Range("E19").Value = '  +0.38%  '
$ws.Range("E20").Value = '  +0.65%  '
$ws.Range("E21").Value = '  +4.21%  '
$ws.Range("D22").Value = "'76.70"
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("D23").Value = "'255.92"
$ws.Range("E23").Value = '  +1.02%  '
$ws.Range("D24").Value = "'1.94"
$ws.Range("E24").Value = '  +22.02%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("E26").Value = '  -1.73%  '
$ws.Range("D27").Value = "'2.48"
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("D28").Value = "'10.68"
$ws.Range("E28").Value = '  +1.18%  '
$ws.Range("D29").Value = "'2.25"
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").Value = "'22.68"
$ws.Range("E30").Value = '  +1.41%  '
$ws.Range("D31").Value = "'175.19"
$ws.Range("E31").Value = '  +1.45%  '
$ws.Range("E32").Value = '  -3.07%  '
$ws.Range("E33").Value = '  +3.03%  '
$ws.Range("D34").Value = "'0.0757"
$ws.Range("E34").Value = '  +6.05%  '
$ws.Range("E35").Value = '  +1.92%  '
$ws.Range("D36").Value = "'5.34"
$ws.Range("E36").Value = '  +3.11%  '
$ws.Range("D37").Value = "'3.73"
$ws.Range("E37").Value = '  -3.79%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").Value = "'2.35"
$ws.Range("E38").Value = '  -2.63%  '
$ws.Range("B39").Value = 'THORChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D39").Value = "'6.33"
$ws.Range("E39").Value = '  -0.98%  '
$ws.Range("D40").Value = "'0.0280"
$ws.Range("E40").Value = '  +4.69%  '
$ws.Range("D41").Value = "'0.110"
$ws.Range("E41").Value = '  +12.30%  '
$ws.Range("D42").Value = "'19.14"
$ws.Range("E42").Value = '  -1.78%  '
$ws.Range("E43").Value = '  +11.91%  '
$ws.Range("D44").Value = "'9.07"
$ws.Range("E44").Value = '  +2.63%  '
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("E46").Value = '  +5.42%  '
$ws.Range("E47").Value = '  +1.50%  '
$ws.Range("D48").Value = "'57.36"
$ws.Range("D49").Value = "'2.48"
$ws.Range("E49").Value = '  +8.67%  '
$ws.Range("E50").Value = '  +0.39%  '
$ws.Range("D51").Value = "'100.04"
$ws.Range("E51").Value = '  +2.88%  '
